$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Semestre ideal:" value (row 9, columns B and C)
$ws.Range("B9").Value = "EA-4,EB-5,EQD-4,EQN-5"
$ws.Range("C9").Value = "EA-4,EB-5,EQD-4,EQN-5"

# Update existing "Requisitos:" entry (row 24) to the new prerequisite text
$ws.Range("B24").Value = "LOB1024 -  Mecânica  (Requisito fraco)`n"
$ws.Range("C24").Value = "LOB1024 -  Mecânica  (Requisito fraco)`n"

# Add a new row (25) with an additional "Requisitos:" entry, copying the
# formatting (styles + row height) from row 24
$ws.Range("B24:C24").Copy()
$ws.Range("B25:C25").PasteSpecial(-4122)

$ws.Range("B25").Value = "LOB1052 -  Cálculo III  (Requisito fraco)`n"
$ws.Range("C25").Value = "LOB1052 -  Cálculo III  (Requisito fraco)`n"

$ws.Rows.Item(25).RowHeight = 30
